$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 - Magnesium chloride unit price: update baseline, lower (was formula), upper (was formula)
$ws.Range("E8").Value = 0.38
$ws.Range("G8").Value = 0.349
$ws.Range("I8").Value = 0.411

# Row 9 - Zinc sulfate unit price: update baseline, lower (was formula), upper (was formula)
$ws.Range("E9").Value = 0.795
$ws.Range("G9").Value = 0.657
$ws.Range("I9").Value = 0.931

# Update the active selection on the sheet (matches recorded view state after edit)
$ws.Activate()
$ws.Range("A8:XFD9").Select()
